$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, so existing column B (LOCATION) shifts to C
$ws.Columns.Item(2).Insert()

# Header row
$ws.Range("A1").Value = "fruits_df1"
$ws.Range("B1").Value = "fruits_df2"
$ws.Range("C1").Value = "LOCATION"

# Copy the header style (bold/centered) from A1 onto the new B1/C1 cells
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 - apple (in both)
$ws.Range("B2").Value = "apple"

# Row 3 - banana (in both)
$ws.Range("B3").Value = "banana"

# Row 4 - cherry (in both)
$ws.Range("B4").Value = "cherry"

# Row 5 - grape (fruits_df2 only) -> move value from A to B, clear A
$ws.Range("B5").Value = $ws.Range("A5").Value()
$ws.Range("A5").Value = ""

# Row 6 - kiwi (fruits_df1 only) -> already in A, clear B (was LOCATION shifted already)
$ws.Range("B6").Value = ""

# Row 7 - mango (fruits_df1 only)
$ws.Range("B7").Value = ""

# Row 8 - peach (fruits_df1 only)
$ws.Range("B8").Value = ""

# Row 9 - pear (fruits_df2 only) -> move value from A to B, clear A
$ws.Range("B9").Value = $ws.Range("A9").Value()
$ws.Range("A9").Value = ""

# Row 10 - watermelon (fruits_df2 only) -> move value from A to B, clear A
$ws.Range("B10").Value = $ws.Range("A10").Value()
$ws.Range("A10").Value = ""
